$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1)
$ws.Range("A1").Value = "MIGRATION DATE"
$ws.Range("B1").Value = "FINANCIAL INSTITUTION NAME"
$ws.Range("C1").Value = "ENTITY ID"
$ws.Range("D1").Value = "ADDRESS"

# C1:D1 should carry the same formatting already present on A1/B1
$ws.Range("A1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)

# Data row (row 2). Force the cells to text first so the date-like value in
# A2 is kept as literal text "2025-10-17" instead of being auto-converted to
# a date serial number, then drop back to the default (unformatted) style.
$ws.Range("A2:D2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-10-17"
$ws.Range("B2").Value = "ZZZ"
$ws.Range("C2").Value = "456CDX009"
$ws.Range("D2").Value = "Anna Nagar"
$ws.Range("A2:D2").ClearFormats()
